# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values; update rows 2-8 with the regenerated values.
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 8
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 7
$ws.Range("G6").Value = 8
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 3
